$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "abcded"
$ws.Range("C2").Value = 4275
$ws.Range("F2").Value = "14:30"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2024-09-03"
